# Add the "11.12" day entry (row 12) to the Adventskalender worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Values are entered in the same order the author's session produced them
# (C, D, E, then B) so the shared-string table ends up with identical
# indices/ordering to the authored workbook.
$ws.Range("C12").Value = "Endlich bin ich wieder Joggen gegangen. Ich habe mir mitten im Nachmittag die Zeit genommen, eine kurze Runde nach draussen zu gehen. Die Aktivierungsphase hat somit wieder gestartet. 😀 Zudem habe ich abends am 10.11 im Internet nach Kochkursen geschaut, die ich buchen könnte, um alleine hinzugehen. Ich möchte schon lange einen Kurs machen, bei dem ich Grundkochskills erlerne. Schlaue Kurse habe ich noch nie wirklich gefunden - bis dahin. Jetzt habe ich mir einen Zugang zur Onlineplattform 7Hauben gekauft. Hier kann ich eben diese Skills erlernen, aber auch spannende Kurse zu gewissen Themen dann machen, wenn es mir geht. Und hier habe ich auch ein Kurs zum Thema `"fermentierte Getränke`" gefunden. Das ist schon lange ein Wunsch von mir, da mich dies mega interessiert und ich die alkoholfreie Getränkebegleitung als sehr reizvolle Aufgabe wahrnehme. In diesem Sinne schauen wir mal, was ich Neues lernen darf."
$ws.Range("D12").Value = "Aktivierung und Neues lernen"
$ws.Range("E12").Value = "https://d2z9jv66wc4wox.cloudfront.net/eyJidWNrZXQiOiI3aC1zdHJhcGkiLCJrZXkiOiI3aGF1YmVuX3Bhc3NfZ3V0c2NoZWluX2RjZDc0OTA3MDMuanBlZyIsImVkaXRzIjp7InJlc2l6ZSI6eyJ3aWR0aCI6NzAwLCJoZWlnaHQiOm51bGwsImZpdCI6ImNvdmVyIn19fQ=="
$ws.Range("B12").Value = "Körperliche Aktivierung & Neues lernen"

# Move the active cell selection to B13, matching the author's next edit location.
$ws.Range("B13").Select()
